$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '23.714.77'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -3.06%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.613.56'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -3.23%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.005'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.15%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '1.004'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +0.15%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '306.29'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -2.00%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3908'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.31%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3802'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -3.10%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '1.005'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +0.16%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.343'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -3.73%  '

$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -6.20%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.08416'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -2.21%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '23.70'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -4.97%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.963'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -4.00%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.00001267'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -3.45%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '7.382'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -4.01%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.611.53'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -3.86%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '93.41'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +0.49%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06914'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -2.12%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '19.82'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -3.83%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.779'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -3.57%  '

$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +0.02%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '13.34'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -4.08%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '23.749.92'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -2.90%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.405'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +1.28%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.796'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +2.16%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '22.04'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -4.56%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '157.46'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -2.02%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '139.11'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -5.54%  '

$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -9.82%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.730'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -6.24%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.472'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -1.71%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.789.82'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -4.07%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.07991'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -3.78%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.9560'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -2.20%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '6.593'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -4.95%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02847'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -5.40%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.2651'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -5.17%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.09123'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -3.20%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '10.33'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +0.44%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '13.22'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -1.68%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.416'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -7.88%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.7424'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -5.46%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '15.80'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -3.32%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.6802'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -3.75%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.428'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -4.10%  '

$ws.Range('B47').NumberFormat = '@'
$ws.Range('B47').Value = 'PancakeSwap'
$ws.Range('C47').NumberFormat = '@'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '4.042'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -3.04%  '

$ws.Range('B48').NumberFormat = '@'
$ws.Range('B48').Value = 'Frax'
$ws.Range('C48').NumberFormat = '@'
$ws.Range('C48').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.004'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +0.18%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.08193'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -4.34%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '132.07'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -3.62%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.241'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -5.62%  '
